$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 64 (hunk 0)
$ws.Range("H64").Value = 3015.1155
$ws.Range("I64").Value = 2799.611
$ws.Range("K64").Value = 2799.611
$ws.Range("M64").Value = -2551.611
# row 67 (hunk 1)
$ws.Range("H67").Value = 3015.1155
$ws.Range("I67").Value = 2799.611
$ws.Range("K67").Value = 2799.611
$ws.Range("M67").Value = -1941.611
# row 76 (hunk 2)
$ws.Range("H76").Value = 15159342
$ws.Range("I76").Value = 9649.9375
$ws.Range("J76").Value = 55558524
$ws.Range("K76").Value = 9649.9375
$ws.Range("L76").Value = 55558524
$ws.Range("M76").Value = -9334.9375
$ws.Range("N76").Value = -55559154
# row 79 (hunk 3)
$ws.Range("H79").Value = 15159342
$ws.Range("I79").Value = 9649.9375
$ws.Range("J79").Value = 55558524
$ws.Range("K79").Value = 9649.9375
$ws.Range("L79").Value = 55558524
$ws.Range("M79").Value = -8557.9375
$ws.Range("N79").Value = -55560708
# row 112 (hunk 4)
$ws.Range("H112").Value = 18977.152
$ws.Range("I112").Value = 533.1667
$ws.Range("J112").Value = 21743.75
$ws.Range("K112").Value = 1599.5001
$ws.Range("L112").Value = 65231.25
$ws.Range("M112").Value = -491.5001
$ws.Range("N112").Value = -67447.25
# row 113 (hunk 5)
$ws.Range("H113").Value = 3615.8235
$ws.Range("I113").Value = 3867.5
$ws.Range("J113").Value = 3538.3845
$ws.Range("K113").Value = 3867.5
$ws.Range("L113").Value = 3538.3845
$ws.Range("M113").Value = -613.5
$ws.Range("N113").Value = -10046.3845
# row 116 (hunk 6)
$ws.Range("H116").Value = 1992.3636
$ws.Range("I116").Value = 1538.75
$ws.Range("J116").Value = 3202
$ws.Range("K116").Value = 1538.75
$ws.Range("L116").Value = 3202
$ws.Range("M116").Value = 1903.25
$ws.Range("N116").Value = -10086
# row 132 (hunk 7)
$ws.Range("H132").Value = 5379756.5
$ws.Range("I132").Value = 5816835
$ws.Range("J132").Value = 9935
$ws.Range("K132").Value = 17450505
$ws.Range("L132").Value = 29805
$ws.Range("M132").Value = -17447975
$ws.Range("N132").Value = -34865
# row 137 (hunk 8)
$ws.Range("H137").Value = 1215.2368
$ws.Range("I137").Value = 730.9167
$ws.Range("J137").Value = 2045.5
$ws.Range("K137").Value = 2192.7501
$ws.Range("L137").Value = 6136.5
$ws.Range("M137").Value = 357.2498999999998
$ws.Range("N137").Value = -11236.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 63 (hunk 9)
$ws.Range("H63").Value = 2479.6
$ws.Range("I63").Value = 2266.3333
$ws.Range("J63").Value = 2799.5
$ws.Range("K63").Value = 2266.3333
$ws.Range("L63").Value = 2799.5
$ws.Range("M63").Value = -1580.3333
$ws.Range("N63").Value = -4171.5
# row 66 (hunk 10)
$ws.Range("H66").Value = 2479.6
$ws.Range("I66").Value = 2266.3333
$ws.Range("J66").Value = 2799.5
$ws.Range("K66").Value = 11331.6665
$ws.Range("L66").Value = 13997.5
$ws.Range("M66").Value = -7899.666499999999
$ws.Range("N66").Value = -20861.5
# row 74 (hunk 11)
$ws.Range("H74").Value = 10870333
$ws.Range("I74").Value = 12500733
$ws.Range("J74").Value = 999.5
$ws.Range("K74").Value = 12500733
$ws.Range("L74").Value = 999.5
$ws.Range("M74").Value = -12499859
$ws.Range("N74").Value = -2747.5
# row 77 (hunk 12)
$ws.Range("H77").Value = 10870333
$ws.Range("I77").Value = 12500733
$ws.Range("J77").Value = 999.5
$ws.Range("K77").Value = 62503665
$ws.Range("L77").Value = 4997.5
$ws.Range("M77").Value = -62499297
$ws.Range("N77").Value = -13733.5
# row 88 (hunk 13)
$ws.Range("H88").Value = 2569.0908
$ws.Range("I88").Value = 3330.111
$ws.Range("J88").Value = 2042.2307
$ws.Range("K88").Value = 3330.111
$ws.Range("L88").Value = 2042.2307
$ws.Range("M88").Value = -2924.111
$ws.Range("N88").Value = -2854.2307
# row 91 (hunk 14)
$ws.Range("H91").Value = 2569.0908
$ws.Range("I91").Value = 3330.111
$ws.Range("J91").Value = 2042.2307
$ws.Range("K91").Value = 3330.111
$ws.Range("L91").Value = 2042.2307
$ws.Range("M91").Value = -1926.111
$ws.Range("N91").Value = -4850.2307
# row 132 (hunk 15)
$ws.Range("H132").Value = 35719610
$ws.Range("I132").Value = 50001480
$ws.Range("K132").Value = 150004440
$ws.Range("M132").Value = -150001910

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 86 (hunk 16)
$ws.Range("H86").Value = 804046.9
$ws.Range("I86").Value = 2803.6924
$ws.Range("J86").Value = 1455057
$ws.Range("K86").Value = 2803.6924
$ws.Range("L86").Value = 1455057
$ws.Range("M86").Value = -1680.6924
$ws.Range("N86").Value = -1457303
# row 89 (hunk 17)
$ws.Range("H89").Value = 804046.9
$ws.Range("I89").Value = 2803.6924
$ws.Range("J89").Value = 1455057
$ws.Range("K89").Value = 14018.462
$ws.Range("L89").Value = 7275285
$ws.Range("M89").Value = -8402.462
$ws.Range("N89").Value = -7286517

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31 (hunk 18)
$ws.Range("H31").Value = 1391.6884
$ws.Range("I31").Value = 983.2632
$ws.Range("J31").Value = 2555.7
$ws.Range("K31").Value = 983.2632
$ws.Range("L31").Value = 2555.7
$ws.Range("M31").Value = -688.2632
$ws.Range("N31").Value = -3145.7
# row 34 (hunk 19)
$ws.Range("H34").Value = 1391.6884
$ws.Range("I34").Value = 983.2632
$ws.Range("J34").Value = 2555.7
$ws.Range("K34").Value = 983.2632
$ws.Range("L34").Value = 2555.7
$ws.Range("M34").Value = -781.2632
$ws.Range("N34").Value = -2959.7
# row 62 (hunk 20)
$ws.Range("H62").Value = 4115.3335
$ws.Range("I62").Value = 2634.6667
$ws.Range("J62").Value = 4608.8887
$ws.Range("K62").Value = 2634.6667
$ws.Range("L62").Value = 4608.8887
$ws.Range("M62").Value = -2010.6667
$ws.Range("N62").Value = -5856.8887
# row 65 (hunk 21)
$ws.Range("H65").Value = 4115.3335
$ws.Range("I65").Value = 2634.6667
$ws.Range("J65").Value = 4608.8887
$ws.Range("K65").Value = 13173.3335
$ws.Range("L65").Value = 23044.4435
$ws.Range("M65").Value = -10053.3335
$ws.Range("N65").Value = -29284.4435

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 4 (hunk 22)
$ws.Range("H4").Value = 333700
$ws.Range("J4").Value = 1000
$ws.Range("L4").Value = 3000
$ws.Range("N4").Value = -3224
# row 81 (hunk 23)
$ws.Range("H81").Value = 10205637
$ws.Range("I81").Value = 14286192
$ws.Range("K81").Value = 42858576
$ws.Range("M81").Value = -42857453
# row 84 (hunk 24)
$ws.Range("H84").Value = 10205637
$ws.Range("I84").Value = 14286192
$ws.Range("K84").Value = 128575728
$ws.Range("M84").Value = -128570112
# row 88 (hunk 25)
$ws.Range("H88").Value = 4490
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 4490
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 13470
$ws.Range("N88").Value = -14326
$ws.Range("M88").ClearContents()
# row 91 (hunk 26)
$ws.Range("H91").Value = 4490
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 4490
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 13470
$ws.Range("N91").Value = -16434
$ws.Range("M91").ClearContents()
# row 122 (hunk 27)
$ws.Range("H122").Value = 14709393
$ws.Range("I122").Value = 55555816
$ws.Range("J122").Value = 4680.92
$ws.Range("K122").Value = 500002344
$ws.Range("L122").Value = 42128.28
$ws.Range("M122").Value = -499999894
$ws.Range("N122").Value = -47028.28
# row 131 (hunk 28)
$ws.Range("H131").Value = 720.4123499999999
$ws.Range("J131").Value = 781.5
$ws.Range("L131").Value = 2344.5
$ws.Range("N131").Value = -12424.5
# row 140 (hunk 29)
$ws.Range("H140").Value = 69446740
$ws.Range("I140").Value = 125001850
$ws.Range("J140").Value = 2849.875
$ws.Range("K140").Value = 375005550
$ws.Range("L140").Value = 8549.625
$ws.Range("M140").Value = -375000370
$ws.Range("N140").Value = -18909.625

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 70 (hunk 30)
$ws.Range("H70").Value = 4622.1113
$ws.Range("I70").Value = 4622.1113
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 4622.1113
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -4352.1113
$ws.Range("N70").ClearContents()
# row 73 (hunk 31)
$ws.Range("H73").Value = 4622.1113
$ws.Range("I73").Value = 4622.1113
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 4622.1113
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -3686.1113
$ws.Range("N73").ClearContents()
# row 80 (hunk 32)
$ws.Range("H80").Value = 7145142
$ws.Range("I80").Value = 2488
$ws.Range("J80").Value = 20001920
$ws.Range("K80").Value = 2488
$ws.Range("L80").Value = 20001920
$ws.Range("M80").Value = -1490
$ws.Range("N80").Value = -20003916
# row 83 (hunk 33)
$ws.Range("H83").Value = 7145142
$ws.Range("I83").Value = 2488
$ws.Range("J83").Value = 20001920
$ws.Range("K83").Value = 12440
$ws.Range("L83").Value = 100009600
$ws.Range("M83").Value = -7448
$ws.Range("N83").Value = -100019584
# row 97 (hunk 34)
$ws.Range("H97").Value = 665.15
$ws.Range("I97").Value = 632.8889
$ws.Range("J97").Value = 955.5
$ws.Range("K97").Value = 632.8889
$ws.Range("L97").Value = 955.5
$ws.Range("M97").Value = -136.8889
$ws.Range("N97").Value = -1947.5
# row 132 (hunk 35)
$ws.Range("H132").Value = 3804.9714
$ws.Range("I132").Value = 746.03845
$ws.Range("J132").Value = 12641.889
$ws.Range("K132").Value = 2238.11535
$ws.Range("L132").Value = 37925.667
$ws.Range("M132").Value = 291.88465
$ws.Range("N132").Value = -42985.667

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 81 (hunk 36)
$ws.Range("H81").Value = 13800
$ws.Range("I81").Value = 13800
$ws.Range("K81").Value = 13800
$ws.Range("M81").Value = -12802
# row 84 (hunk 37)
$ws.Range("H84").Value = 13800
$ws.Range("I84").Value = 13800
$ws.Range("K84").Value = 41400
$ws.Range("M84").Value = -36408

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 81 (hunk 38)
$ws.Range("H81").Value = 27778376
$ws.Range("I81").Value = 38462050
$ws.Range("J81").Value = 830.8
$ws.Range("K81").Value = 76924100
$ws.Range("L81").Value = 1661.6
$ws.Range("M81").Value = -76923039
$ws.Range("N81").Value = -3783.6
# row 84 (hunk 39)
$ws.Range("H84").Value = 27778376
$ws.Range("I84").Value = 38462050
$ws.Range("J84").Value = 830.8
$ws.Range("K84").Value = 384620500
$ws.Range("L84").Value = 8308
$ws.Range("M84").Value = -384615196
$ws.Range("N84").Value = -18916

